# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt"
# (Uva / Red Globe) at row 345, pushing the existing rows 345-359 down to
# 346-360 (dimension grows from A1:T359 to A1:T360).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(345).Insert()

$ws.Cells.Item(345, 1).Value = 4
$ws.Cells.Item(345, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(345, 3).Value = "Los Lagos"
$ws.Cells.Item(345, 4).Value = 45041
$ws.Cells.Item(345, 5).Value = 10
$ws.Cells.Item(345, 6).Value = "Fruta"
$ws.Cells.Item(345, 7).Value = 100109
$ws.Cells.Item(345, 8).Value = "Uva"
$ws.Cells.Item(345, 9).Value = 100109001
$ws.Cells.Item(345, 10).Value = "Uva"
$ws.Cells.Item(345, 11).Value = "Red Globe"
$ws.Cells.Item(345, 12).Value = "Primera"
$ws.Cells.Item(345, 13).Value = 200
$ws.Cells.Item(345, 14).Value = 14000
$ws.Cells.Item(345, 15).Value = 15000
$ws.Cells.Item(345, 16).Value = 14500
$ws.Cells.Item(345, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(345, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(345, 19).Value = 806
$ws.Cells.Item(345, 20).Value = 18
